$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.434.99"
$ws.Range("E2").Value = "  +0.75%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.303.29"
$ws.Range("E3").Value = "  -0.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.18"
$ws.Range("E5").Value = "  +0.90%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.14"
$ws.Range("E6").Value = "  -2.26%  "

# Row 7
$ws.Range("E7").Value = "  +0.42%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("E9").Value = "  +0.22%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.15"
$ws.Range("E10").Value = "  -0.90%  "

# Row 11
$ws.Range("E11").Value = "  -0.68%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.38"
$ws.Range("E12").Value = "  +0.87%  "

# Row 13
$ws.Range("E13").Value = "  +0.56%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.970"
$ws.Range("E14").Value = "  -1.52%  "

# Row 15
$ws.Range("E15").Value = "  -2.14%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.649.94"
$ws.Range("E16").Value = "  -0.60%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.309.85"
$ws.Range("E17").Value = "  -0.39%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.511.17"
$ws.Range("E18").Value = "  +0.64%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.49"
$ws.Range("E19").Value = "  -2.84%  "

# Row 20
$ws.Range("E20").Value = "  +0.98%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.31"
$ws.Range("E21").Value = "  -2.08%  "

# Row 22
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "277.02"
$ws.Range("E22").Value = "  +6.13%  "

# Row 23
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.53"
$ws.Range("E23").Value = "  +0.77%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.04"
$ws.Range("E24").Value = "  +18.32%  "

# Row 25
$ws.Range("E25").Value = "  -1.06%  "

# Row 26
$ws.Range("E26").Value = "  -0.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.84"
$ws.Range("E27").Value = "  -1.74%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.35"
$ws.Range("E28").Value = "  +3.28%  "

# Row 29
$ws.Range("E29").Value = "  -0.24%  "

# Row 30
$ws.Range("E30").Value = "  -0.66%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.58"
$ws.Range("E31").Value = "  +0.27%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0875"
$ws.Range("E32").Value = "  -3.04%  "

# Row 33
$ws.Range("E33").Value = "  +0.35%  "

# Row 34
$ws.Range("E34").Value = "  +4.73%  "

# Row 35
$ws.Range("E35").Value = "  -10.52%  "

# Row 36
$ws.Range("E36").Value = "  -3.08%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0370"
$ws.Range("E37").Value = "  +4.47%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.60"
$ws.Range("E38").Value = "  +1.22%  "

# Row 39
$ws.Range("E39").Value = "  +2.52%  "

# Row 40
$ws.Range("E40").Value = "  -1.23%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.50"
$ws.Range("E41").Value = "  +1.89%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.04"
$ws.Range("E42").Value = "  -1.74%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.227"
$ws.Range("E43").Value = "  -1.26%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.53"
$ws.Range("E45").Value = "  -5.13%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "82.76"
$ws.Range("E46").Value = "  +10.83%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.05"
$ws.Range("E47").Value = "  -3.03%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "113.07"
$ws.Range("E48").Value = "  +0.67%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.591.64"
$ws.Range("E50").Value = "  +2.54%  "

# Row 51
$ws.Range("E51").Value = "  -5.65%  "
